$d = $word.ActiveDocument

# Replace the placeholder "[NEED TIME HERE]," with the concrete "180 days,"
# text (and implicitly drop the yellow highlight formatting that was only
# applied to the placeholder run) in both occurrences of this reminder
# sentence in the document.
$d.Content.Find.Execute("If it has been over [NEED TIME HERE], ", $true, $false, $false, $false, $false,
                         $false, 1, $false, "If it has been over 180 days, ", 2)
